# Auto-committed change: update the "形態" (data type) for CreateDate and
# LastUpdate fields on the DBD sheet from DATE to TIMESTAMP.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")

$ws.Range("D14").Value = "TIMESTAMP"
$ws.Range("D16").Value = "TIMESTAMP"

# Match the author's final selection/active cell on the sheet.
$ws.Range("D16").Select() | Out-Null
